# Updated cryptos list (Price / Volume(1h) refresh, plus a couple of
# ranking swaps) as produced by the scheduled GitHub Actions scraper run.
#
# Price values in column D look like plain numbers (e.g. "1.000", "323.12")
# but must stay as literal text (matching the original "29.324.77" style
# thousands-grouped / zero-padded strings), so NumberFormat is forced to
# "@" (Text) right before assigning them. ClearFormats() afterwards drops
# the now-unneeded explicit style index again so the cells end up with no
# "s" attribute, same as every other untouched data cell in that column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '29.291.07'
$ws.Range('E2').Value = '  +0.92%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.905.81'
$ws.Range('E3').Value = '  +1.00%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '323.12'
$ws.Range('E5').Value = '  -2.35%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4724'
$ws.Range('E7').Value = '  +2.62%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4030'
$ws.Range('E8').Value = '  -1.09%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '47.73'
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('E10').Value = '  +0.46%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.9932'
$ws.Range('E11').Value = '  +0.21%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '22.66'
$ws.Range('E12').Value = '  +4.37%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.891.03'
$ws.Range('E13').Value = '  +0.44%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.854'
$ws.Range('E14').Value = '  -1.02%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.042'
$ws.Range('E15').Value = '  -0.44%  '
$ws.Range('E16').Value = '  +1.02%  '
$ws.Range('E17').Value = '  +0.06%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.06605'
$ws.Range('E18').Value = '  +0.44%  '
$ws.Range('E19').Value = '  -0.51%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.52'
$ws.Range('E20').Value = '  +0.37%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.002'
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '29.309.95'
$ws.Range('E22').Value = '  +0.98%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.501'
$ws.Range('E23').Value = '  +1.39%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.43'
$ws.Range('E24').Value = '  -0.20%  '
$ws.Range('E25').Value = '  -0.27%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.129.05'
$ws.Range('E26').Value = '  +1.04%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '153.90'
$ws.Range('E27').Value = '  -1.77%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '19.73'
$ws.Range('E28').Value = '  +0.82%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.030'
$ws.Range('E29').Value = '  +10.09%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.092'
$ws.Range('E30').Value = '  +0.26%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '117.84'
$ws.Range('E31').Value = '  +0.28%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.072'
$ws.Range('E32').Value = '  +6.40%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.09502'
$ws.Range('E33').Value = '  +1.98%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.412'
$ws.Range('E34').Value = '  +0.32%  '
$ws.Range('E35').Value = '  -1.42%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.354'
$ws.Range('E36').Value = '  +1.48%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.06060'
$ws.Range('E37').Value = '  -0.17%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.02245'
$ws.Range('E38').Value = '  +0.75%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.171'
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '8.084'
$ws.Range('E40').Value = '  -2.60%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.5802'
$ws.Range('E41').Value = '  +0.34%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.496'
$ws.Range('E42').Value = '  +9.40%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.1831'
$ws.Range('E43').Value = '  +0.30%  '
$ws.Range('E44').Value = '  -0.75%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.07821'
$ws.Range('E45').Value = '  +4.11%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.278'
$ws.Range('E46').Value = '  +0.92%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '12.18'
$ws.Range('E47').Value = '  +1.21%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.5479'
$ws.Range('E48').Value = '  +0.51%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.893'
$ws.Range('E49').Value = '  -0.48%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '113.14'
$ws.Range('E50').Value = '  +1.83%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '43.95'
$ws.Range('E51').Value = '  -3.46%  '

# Restore default (no explicit) style on the Price column data cells that were
# forced to Text format above, so their cell style matches the original (no s attribute).
$ws.Range('D2:D51').ClearFormats()
